# New crime data collected
# Updates the weekly NYPD CompStat report (1st Precinct) with the new
# reporting week's figures: header "Volume/Number" + reporting date
# range, and the crime-complaint statistics table (rows 15-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text -----------------------------------------------------
# "Volume 32   Number  8"  ->  "Volume 32   Number  9"
$ws.Range("A8").Value = "Volume 32   Number  9"
# "Report Covering the Week  2/17/2025  Through  2/23/2025"
#   -> "...2/24/2025  Through  3/2/2025"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Row 15 : Rape -----------------------------------------------------
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 0

# --- Row 16 : Robbery ----------------------------------------------------
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 27.272727272727
$ws.Range("I16").Value = 21
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = -16
$ws.Range("L16").Value = -8.695652173913
$ws.Range("M16").Value = 133.333333333333
$ws.Range("N16").Value = -87.037037037037

# --- Row 17 : Fel. Assault ----------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = -36
$ws.Range("L17").Value = -5.882352941176
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = -38.461538461538

# --- Row 18 : Burglary ---------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 44
$ws.Range("K18").Value = -40.909090909090
$ws.Range("L18").Value = -7.142857142857
$ws.Range("M18").Value = -33.333333333333
$ws.Range("N18").Value = -82.312925170068

# --- Row 19 : Gr. Larceny -------------------------------------------------
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = 15.384615384615
$ws.Range("F19").Value = 106
$ws.Range("G19").Value = 85
$ws.Range("H19").Value = 24.705882352941
$ws.Range("I19").Value = 199
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = 17.058823529411
$ws.Range("L19").Value = 11.173184357541
$ws.Range("M19").Value = 11.797752808988
$ws.Range("N19").Value = -64.014466546112

# --- Row 20 : G.L.A. -------------------------------------------------------
$ws.Range("N20").Value = -97.419354838709

# --- Row 21 : TOTAL ---------------------------------------------------------
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = -4.651162790697
$ws.Range("F21").Value = 150
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = 6.382978723404
$ws.Range("I21").Value = 270
$ws.Range("J21").Value = 272
$ws.Range("K21").Value = -0.735294117647
$ws.Range("L21").Value = 5.46875
$ws.Range("M21").Value = 13.445378151260
$ws.Range("N21").Value = -74.310180780209

# --- Row 22 : Transit -------------------------------------------------------
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 10
$ws.Range("H22").Value = 42.857142857142
$ws.Range("I22").Value = 14
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = 7.692307692307
$ws.Range("L22").Value = -17.647058823529
$ws.Range("M22").Value = 0

# --- Row 24 : Petit Larceny --------------------------------------------------
$ws.Range("C24").Value = 74
$ws.Range("D24").Value = 94
$ws.Range("E24").Value = -21.276595744680
$ws.Range("F24").Value = 289
$ws.Range("G24").Value = 319
$ws.Range("H24").Value = -9.404388714733
$ws.Range("I24").Value = 596
$ws.Range("J24").Value = 677
$ws.Range("K24").Value = -11.964549483013
$ws.Range("L24").Value = -5.696202531645
$ws.Range("M24").Value = 154.700854700855

# --- Row 25 : Retail Theft ----------------------------------------------------
$ws.Range("C25").Value = 83
$ws.Range("D25").Value = 89
$ws.Range("E25").Value = -6.741573033707
$ws.Range("F25").Value = 301
$ws.Range("G25").Value = 318
$ws.Range("H25").Value = -5.345911949685
$ws.Range("I25").Value = 599
$ws.Range("J25").Value = 675
$ws.Range("K25").Value = -11.259259259259
$ws.Range("L25").Value = -5.817610062893

# --- Row 26 : Misd. Assault ---------------------------------------------------
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = -29.729729729729
$ws.Range("I26").Value = 69
$ws.Range("J26").Value = 58
$ws.Range("K26").Value = 18.965517241379
$ws.Range("L26").Value = 43.75
$ws.Range("M26").Value = 97.142857142857

# --- Row 27 : UCR Rape* -----------------------------------------------------
# C27 used to be the blank-placeholder text "0"; it now carries a real
# number, so give it the same numeric style/format used by its siblings
# (e.g. F27/G27) rather than the inherited text style.
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 150

# --- Row 28 : Other Sex Crimes -----------------------------------------------
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -14.285714285714
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = -46.666666666666

# --- Row 31 : Hate Crimes ---------------------------------------------------
$ws.Range("I31").Value = 3
$ws.Range("L31").Value = -40
